# wyniki/1.xlsx update:
#  - E2:E12 scores become real numbers instead of numeric-looking text
#  - three new columns (G: Nazwa Pliku i rozdzial / H: Wynik / I: Maksymalny Wynik)
#  - four new data rows (13-16) for the new audio-file based entries
#  - a couple of "present but empty" cells that existed in the source export

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. Fix E2:E12 so the score is stored as a number, not text ----
$scores = @{
    2  = 2
    3  = 0
    4  = 5
    5  = 2
    6  = 1
    7  = 2
    8  = 2
    9  = 1
    10 = 5
    11 = 3
    12 = 2
}
foreach ($row in $scores.Keys) {
    $ws.Cells.Item($row, 5).Value = $scores[$row]
}

# ---- 2. New header cells G1:I1, styled like the existing header row ----
$ws.Range("G1").Value = "Nazwa Pliku i rozdział"
$ws.Range("H1").Value = "Wynik"
$ws.Range("I1").Value = "Maksymalny Wynik"

$headerRng = $ws.Range("G1:I1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

# ---- 3. C3 and G2:I12 existed as empty (but present/typed) cells in the
#         original export - force them to exist without visible content.
#         A lone "'" quote-prefix stores as an empty text cell (no visible
#         content, type = text) instead of leaving the cell absent. ----
$ws.Range("C3").Value = "'"
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 7).Value = "'"
    $ws.Cells.Item($row, 8).Value = "'"
    $ws.Cells.Item($row, 9).Value = "'"
}

# ---- 4. Four new rows of data (13-16) ----
$newRows = @(
    @{ Row = 13; A = 12; B = " Piękne,  przykładem"; C = "przykładem."; D = " Piękne,  przykładem"; G = "TAYLOR__Mechanika16k.wav"; H = "1/2"; I = 2 },
    @{ Row = 14; A = 13; B = " Piękne,  przykładem"; C = "przykładem."; D = " Piękne,  przykładem"; G = "TAYLOR__Mechanika16k.wav"; H = "1/2"; I = 2 },
    @{ Row = 15; A = 14; B = " się,  ono,  odepchnąć,,  lub,  który,  mógłby,  ciało"; C = "Ono"; D = " ciało,  mógłby,  ono,  się,  który,  lub,  odepchnąć,"; G = "TAYLOR__Mechanika16k.wav"; H = "1/7"; I = 7 },
    @{ Row = 16; A = 15; B = " się,  ono,  odepchnąć,,  lub,  który,  mógłby,  ciało"; C = "Ono"; D = " ciało,  mógłby,  ono,  się,  który,  lub,  odepchnąć,"; G = "TAYLOR__Mechanika16k.wav"; H = "1/7"; I = 7 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    # E and F stay empty for these rows, but were present/typed cells in the export
    $ws.Cells.Item($row, 5).Value = "'"
    $ws.Cells.Item($row, 6).Value = "'"
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
}

Write-Output "edit applied"
